# Applies the weekly Fruta/Hortaliza data refresh for the Zanahoria - Terminal La Palmera de La Serena subsheet.
# Two new daily records are spliced in (pushing subsequent rows down) and the series is
# extended with three additional rows at the end, per the updated source extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D, J, K, L, M, O, P for target rows 151-223 (A,B,C,E,F,G,H,I,N,Q,R stay constant).
$rows = @(
    @(151, 44490, 600, 7000, 7500, 7250, "Provincia del Elquí", 362),
    @(152, 44396, 700, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(153, 44399, 720, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(154, 44357, 660, 4800, 5000, 4900, "Provincia del Elquí", 245),
    @(155, 44329, 680, 4800, 5000, 4900, "Provincia del Elquí", 245),
    @(156, 44229, 540, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(157, 44376, 600, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(158, 44316, 800, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(159, 44453, 700, 4500, 5000, 4750, "Provincia del Elquí", 238),
    @(160, 44466, 760, 6000, 7000, 6500, "Provincia del Elquí", 325),
    @(161, 44242, 680, 5500, 6000, 5750, "Provincia del Elquí", 288),
    @(162, 44351, 800, 4800, 5000, 4900, "Provincia del Elquí", 245),
    @(163, 44279, 800, 6500, 7000, 6750, "Chillán", 338),
    @(164, 44279, 600, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(165, 44280, 680, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(166, 44389, 760, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(167, 44442, 720, 4500, 5000, 4750, "Provincia del Elquí", 238),
    @(168, 44476, 560, 6000, 7000, 6500, "Provincia del Elquí", 325),
    @(169, 44372, 700, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(170, 44201, 560, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(171, 44482, 800, 6000, 7000, 6500, "Provincia del Elquí", 325),
    @(172, 44293, 800, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(173, 44204, 740, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(174, 44277, 700, 6500, 7000, 6750, "Chillán", 338),
    @(175, 44218, 760, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(176, 44221, 520, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(177, 44166, 560, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(178, 44477, 800, 6000, 7000, 6500, "Provincia del Elquí", 325),
    @(179, 44292, 600, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(180, 44508, 600, 6500, 7000, 6750, "Provincia del Elquí", 338),
    @(181, 44333, 660, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(182, 44356, 800, 4800, 5000, 4900, "Provincia del Elquí", 245),
    @(183, 44323, 800, 4800, 5000, 4900, "Provincia del Elquí", 245),
    @(184, 44306, 600, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(185, 44211, 740, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(186, 44469, 600, 6000, 7000, 6500, "Provincia del Elquí", 325),
    @(187, 44215, 500, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(188, 44407, 720, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(189, 44186, 700, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(190, 44504, 600, 6500, 7000, 6750, "Provincia del Elquí", 338),
    @(191, 44384, 800, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(192, 44246, 800, 5500, 6000, 5750, "Provincia del Elquí", 288),
    @(193, 44505, 760, 6500, 7000, 6750, "Provincia del Elquí", 338),
    @(194, 44487, 600, 6500, 7000, 6750, "Provincia del Elquí", 338),
    @(195, 44425, 660, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(196, 44343, 700, 4800, 5000, 4900, "Provincia del Elquí", 245),
    @(197, 44370, 800, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(198, 44449, 700, 4500, 5000, 4750, "Provincia del Elquí", 238),
    @(199, 44168, 600, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(200, 44175, 600, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(201, 44392, 700, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(202, 44286, 800, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(203, 44473, 600, 6000, 7000, 6500, "Provincia del Elquí", 325),
    @(204, 44400, 720, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(205, 44484, 760, 6500, 7000, 6750, "Provincia del Elquí", 338),
    @(206, 44181, 400, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(207, 44494, 600, 6800, 7000, 6900, "Provincia del Elquí", 345),
    @(208, 44342, 800, 4800, 5000, 4900, "Provincia del Elquí", 245),
    @(209, 44328, 800, 4800, 5000, 4900, "Provincia del Elquí", 245),
    @(210, 44301, 700, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(211, 44330, 800, 4800, 5000, 4900, "Provincia del Elquí", 245),
    @(212, 44270, 600, 5500, 6000, 5750, "Provincia del Elquí", 288),
    @(213, 44295, 800, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(214, 44217, 600, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(215, 44421, 700, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(216, 44383, 600, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(217, 44244, 800, 5500, 6000, 5750, "Provincia del Elquí", 288),
    @(218, 44307, 800, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(219, 44273, 700, 6000, 6500, 6250, "Provincia del Elquí", 312),
    @(220, 44433, 800, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(221, 44302, 800, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(222, 44179, 760, 5000, 5500, 5250, "Provincia del Elquí", 262),
    @(223, 44491, 760, 7000, 7500, 7250, "Provincia del Elquí", 362)
)

$constB = "Terminal La Palmera de La Serena"
$constC = "Coquimbo"
$constN = "$/saco 20 kilos"
$constG = "Zanahoria"
$constH = "Sin especificar"
$constI = "Primera"
$constR = "Hortaliza"

foreach ($r in $rows) {
    $rowNum = $r[0]
    if ($rowNum -gt 220) {
        $ws.Cells.Item($rowNum, 1).Value = 8
        $ws.Cells.Item($rowNum, 2).Value = $constB
        $ws.Cells.Item($rowNum, 3).Value = $constC
        $ws.Cells.Item($rowNum, 5).Value = 4
        $ws.Cells.Item($rowNum, 6).Value = 100114013
        $ws.Cells.Item($rowNum, 7).Value = $constG
        $ws.Cells.Item($rowNum, 8).Value = $constH
        $ws.Cells.Item($rowNum, 9).Value = $constI
        $ws.Cells.Item($rowNum, 14).Value = $constN
        $ws.Cells.Item($rowNum, 17).Value = 20
        $ws.Cells.Item($rowNum, 18).Value = $constR
    }
    $ws.Cells.Item($rowNum, 4).Value = $r[1]
    $ws.Cells.Item($rowNum, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($rowNum, 10).Value = $r[2]
    $ws.Cells.Item($rowNum, 11).Value = $r[3]
    $ws.Cells.Item($rowNum, 12).Value = $r[4]
    $ws.Cells.Item($rowNum, 13).Value = $r[5]
    $ws.Cells.Item($rowNum, 15).Value = $r[6]
    $ws.Cells.Item($rowNum, 16).Value = $r[7]
}

Write-Output "Updated rows 151-223"
